$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.062.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.367.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.64%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.37'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.95'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.357.80'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.638'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.94'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.892.54'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.25'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.346.32'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '64.830.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.84'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.996'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '450.46'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.92'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.70'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.72'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.75%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.77'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.33%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.87'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.68'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.58'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '63.17'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.46'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '575.74'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.11%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.62'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.69%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.71'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0744'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.56%  '

$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.372'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.084.47'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0418'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.77'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.21'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.45%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.73'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.28'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.32%  '
